$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the rule table: C10 (the "min" bound for rule R40) changes from 18 to 1
$ws.Range("C10").Value = 1
